$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $ok = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: replacement failed for:" $find
    }
}

# --- Bullet: cost estimate paragraph (full rewrite) ---
$old6 = "The cost estimate of implementing this law for the next 5 year period is `$500,000. This might seem strange but in fact the DHS has been already conducting these coordination activities and this law is simply codifying what the DHS has already been responsible for thus giving them the complete authority to do so. Over the next 5 years there will actually be no new operating requirements from S. 2520, but looking onward I suspect that this can have two completely different outcomes on Cybersecurity."
$new6 = "The cost estimate for implementing S. 2520 for the next 5 year period is `$500,000. This might seem strange but in fact the DHS has been already conducting these coordination activities mentioned. This law is simply codifying what the DHS has already been responsible for, thus giving them the rekongized authority to do so now. Because of this fact, the next 5 year period after S. 2520 comes into effect there will be actually no new operating requirements required. Looking onward I suspect that this law can have many different outcomes on Cybersecurity and penetration testing. I think that S. 2520 is laying the foundation for further legislation; specifically to increase the efforts to minimize the impact from cybercriminals domestic and abroad. Our government has been overdue on increasing the security of their systems especially due to recent tensions."
Replace-Text $old6 $new6

# --- Bullet: S. 2520 section F / elections paragraph (full rewrite) ---
$old7 = "In S. 2520 section F there is heavy emphasis on promoting the security of elections through collaboration between chief officals of all government entities to enforce policies and procedures related to the security of election systems. I think that in the long term this law is setting the foundation of locking down election systems and treating election systems as critical systems themselves using the same security design and procedures that would give them the fault tolerance of computers of the likes of that power our infrustructure and send send objects to orbit. This is a problem that has caused divide and heart-ache between the American people."
$new7 = "In S. 2520 section F, there is heavy emphasis on promoting the security of elections. This will be done through collaboration between chief officals from all government entities to enforce policies and procedures related to the security of election systems. I think that in the long term this law is setting the foundation of locking down election systems and treating election systems as critical systems themselves, using the same security designs and procedures that would give them reliabiltiy and fault tolerance the likes of the powergrid system and computers that send objects out of our orbit. This is a problem that has caused divide and heart- ache between the American people."
Replace-Text $old7 $new7

# --- Bullet: jobs / education paragraph - small fixes ---
Replace-Text "government and DHS, they will create new standardized education to fill these jobs. Further" "government and the DHS, they will create new standardized education to fill these jobs. Further"

Replace-Text "institution at the community college and upperclass high school level; Funnelling students into roles that help the government at the local levels." "institution at the community college and upperclass high school level can prove to be beneficial; Funnelling these students into roles that help the government at the local levels."

# --- Bullet: black-hat hackers paragraph - add trailing sentence ---
Replace-Text "hosts to government entities." "hosts to government entities. This is a trend that I see to continue for the noticable future."
